# Generate Report for Handoff
# Adds a new "e5b43235-1377-40d3-9218-9b140b175c59" file row (row 3) to the
# Overview, zh-cn and de-de worksheets, mirroring the existing
# "dbd5b577-889d-4181-8830-5ac3c47dbb6b" row.

$wb = $excel.ActiveWorkbook

$guid = "e5b43235-1377-40d3-9218-9b140b175c59"
$srcHash = "5a55b13256be2edf585ed722c7b224ce67598ba8"

$mdDisplay = "$guid.md"
$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/62e8f39d25da217ae0a7f283b0b44347c2f81f13/e2e/$guid.md"

$zhXlfName = "$guid.$srcHash.zh-cn.xlf"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b13304740baf6d87abd1214353266aea59ad923/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"

$deXlfName = "$guid.$srcHash.de-de.xlf"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4d283fc692d98e2da6e39997d28bb50513c966ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-59-17 10:59:00"

$ws1.Range("A3").Value = $mdDisplay
$ws1.Hyperlinks.Add($ws1.Range("A3"), $mdTarget, "", "", $mdDisplay)
$ws1.Range("A3").Font.Underline = $true
$ws1.Range("A3").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("E3").Value = "2016-03-17 10:58:57"
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("I3").Value = "Include"

$ws2.Range("A3").Value = $mdDisplay
$ws2.Hyperlinks.Add($ws2.Range("A3"), $mdTarget, "", "", $mdDisplay)
$ws2.Range("A3").Font.Underline = $true
$ws2.Range("A3").Font.Color = $hyperlinkColor

$ws2.Range("B3").Value = ".md"
$ws2.Hyperlinks.Add($ws2.Range("B3"), $mdTarget, "", "", ".md")
$ws2.Range("B3").Font.Underline = $true
$ws2.Range("B3").Font.Color = $hyperlinkColor

$ws2.Range("D3").Value = $zhXlfName
$ws2.Hyperlinks.Add($ws2.Range("D3"), $zhXlfTarget, "", "", $zhXlfName)
$ws2.Range("D3").Font.Underline = $true
$ws2.Range("D3").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("E3").Value = "2016-03-17 10:59:00"
$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("I3").Value = "Include"

$ws3.Range("A3").Value = $mdDisplay
$ws3.Hyperlinks.Add($ws3.Range("A3"), $mdTarget, "", "", $mdDisplay)
$ws3.Range("A3").Font.Underline = $true
$ws3.Range("A3").Font.Color = $hyperlinkColor

$ws3.Range("B3").Value = ".md"
$ws3.Hyperlinks.Add($ws3.Range("B3"), $mdTarget, "", "", ".md")
$ws3.Range("B3").Font.Underline = $true
$ws3.Range("B3").Font.Color = $hyperlinkColor

$ws3.Range("D3").Value = $deXlfName
$ws3.Hyperlinks.Add($ws3.Range("D3"), $deXlfTarget, "", "", $deXlfName)
$ws3.Range("D3").Font.Underline = $true
$ws3.Range("D3").Font.Color = $hyperlinkColor

Write-Output "Report row added for $guid"
